# Scheduled-runner refresh: update market-board derived profit columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H:N)
# for the affected leve rows across each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2408.7778
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 2659.875
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 7979.625
$ws.Range("M17").Value = -1032
$ws.Range("N17").Value = -8315.625

$ws.Range("H98").Value = 7100.7
$ws.Range("I98").Value = 7791
$ws.Range("J98").Value = 888
$ws.Range("K98").Value = 7791
$ws.Range("L98").Value = 888
$ws.Range("M98").Value = -6293
$ws.Range("N98").Value = -3884

$ws.Range("H116").Value = 3055.7334
$ws.Range("I116").Value = 2395
$ws.Range("J116").Value = 4377.2
$ws.Range("K116").Value = 2395
$ws.Range("L116").Value = 4377.2
$ws.Range("M116").Value = 1047
$ws.Range("N116").Value = -11261.2

$ws.Range("H122").Value = 7100.7
$ws.Range("I122").Value = 7791
$ws.Range("J122").Value = 888
$ws.Range("K122").Value = 23373
$ws.Range("L122").Value = 2664
$ws.Range("M122").Value = -20923
$ws.Range("N122").Value = -7564

$ws.Range("H138").Value = 454605.28
$ws.Range("I138").Value = 918.5263
$ws.Range("J138").Value = 566553.9399999999
$ws.Range("K138").Value = 2755.5789
$ws.Range("L138").Value = 1699661.82
$ws.Range("M138").Value = 2384.4211
$ws.Range("N138").Value = -1709941.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 43479340
$ws.Range("I61").Value = 52632500
$ws.Range("K61").Value = 52632500
$ws.Range("M61").Value = -52632288

$ws.Range("H97").Value = 572.4286
$ws.Range("I97").Value = 543.6667
$ws.Range("J97").Value = 745
$ws.Range("K97").Value = 543.6667
$ws.Range("L97").Value = 745
$ws.Range("M97").Value = -47.66669999999999
$ws.Range("N97").Value = -1737

$ws.Range("H122").Value = 1569.9656
$ws.Range("I122").Value = 1359.5264
$ws.Range("J122").Value = 1969.8
$ws.Range("K122").Value = 4078.5792
$ws.Range("L122").Value = 5909.4
$ws.Range("M122").Value = -1628.5792
$ws.Range("N122").Value = -10809.4

$ws.Range("H136").Value = 43479340
$ws.Range("I136").Value = 52632500
$ws.Range("K136").Value = 157897500
$ws.Range("M136").Value = -157894950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 100001320
$ws.Range("I105").Value = 111112450
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 111112450
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = -111110703
$ws.Range("N105").Value = -4694

$ws.Range("H107").Value = 1294.4706
$ws.Range("I107").Value = 898.5714
$ws.Range("J107").Value = 1571.6
$ws.Range("K107").Value = 898.5714
$ws.Range("L107").Value = 1571.6
$ws.Range("M107").Value = 1021.4286
$ws.Range("N107").Value = -5411.6

$ws.Range("H134").Value = 6872.2104
$ws.Range("I134").Value = 1111.5333
$ws.Range("K134").Value = 3334.5999
$ws.Range("M134").Value = -799.5999000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 611.4
$ws.Range("J10").Value = 500
$ws.Range("L10").Value = 500
$ws.Range("N10").Value = -778

$ws.Range("H23").Value = 2227.25
$ws.Range("I23").Value = 2636.3333
$ws.Range("J23").Value = 1000
$ws.Range("K23").Value = 2636.3333
$ws.Range("L23").Value = 1000
$ws.Range("M23").Value = -2396.3333
$ws.Range("N23").Value = -1480

$ws.Range("H27").Value = 2227.25
$ws.Range("I27").Value = 2636.3333
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 2636.3333
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -2444.3333
$ws.Range("N27").Value = -1384

$ws.Range("H31").Value = 1723.0667
$ws.Range("I31").Value = 2047
$ws.Range("J31").Value = 1439.625
$ws.Range("K31").Value = 2047
$ws.Range("L31").Value = 1439.625
$ws.Range("M31").Value = -1752
$ws.Range("N31").Value = -2029.625

$ws.Range("H34").Value = 1723.0667
$ws.Range("I34").Value = 2047
$ws.Range("J34").Value = 1439.625
$ws.Range("K34").Value = 2047
$ws.Range("L34").Value = 1439.625
$ws.Range("M34").Value = -1845
$ws.Range("N34").Value = -1843.625

$ws.Range("H58").Value = 1588.3704
$ws.Range("I58").Value = 1336.4445
$ws.Range("K58").Value = 1336.4445
$ws.Range("M58").Value = -1133.4445

$ws.Range("H62").Value = 4083983.2
$ws.Range("I62").Value = 2381.7778
$ws.Range("K62").Value = 2381.7778
$ws.Range("M62").Value = -1757.7778

$ws.Range("H65").Value = 4083983.2
$ws.Range("I65").Value = 2381.7778
$ws.Range("K65").Value = 11908.889
$ws.Range("M65").Value = -8788.888999999999

$ws.Range("H95").Value = 13146.5
$ws.Range("J95").Value = 13146.5
$ws.Range("L95").Value = 13146.5
$ws.Range("N95").Value = -18638.5

$ws.Range("H105").Value = 758.5
$ws.Range("I105").Value = 758.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 758.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 988.5
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 470.75
$ws.Range("I107").Value = 368.75
$ws.Range("J107").Value = 674.75
$ws.Range("K107").Value = 368.75
$ws.Range("L107").Value = 674.75
$ws.Range("M107").Value = 1551.25
$ws.Range("N107").Value = -4514.75

$ws.Range("H136").Value = 1588.3704
$ws.Range("I136").Value = 1336.4445
$ws.Range("K136").Value = 4009.3335
$ws.Range("M136").Value = -1459.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6266.75
$ws.Range("I56").Value = 6266.75
$ws.Range("K56").Value = 6266.75
$ws.Range("M56").Value = -5736.75

$ws.Range("H108").Value = 2117.0557
$ws.Range("I108").Value = 482.33334
$ws.Range("J108").Value = 2444
$ws.Range("K108").Value = 1447.00002
$ws.Range("L108").Value = 7332
$ws.Range("M108").Value = 1432.99998
$ws.Range("N108").Value = -13092

$ws.Range("H122").Value = 1190.4117
$ws.Range("J122").Value = 1225.8
$ws.Range("L122").Value = 11032.2
$ws.Range("N122").Value = -15932.2

$ws.Range("H131").Value = 71431610
$ws.Range("J131").Value = 3784.818
$ws.Range("L131").Value = 11354.454
$ws.Range("N131").Value = -21434.454

$ws.Range("H132").Value = 956.8095
$ws.Range("I132").Value = 933.9375
$ws.Range("J132").Value = 1030
$ws.Range("K132").Value = 8405.4375
$ws.Range("L132").Value = 9270
$ws.Range("M132").Value = -5875.4375
$ws.Range("N132").Value = -14330

$ws.Range("H139").Value = 1690.7949
$ws.Range("I139").Value = 1686.32
$ws.Range("J139").Value = 1698.7858
$ws.Range("K139").Value = 5058.96
$ws.Range("L139").Value = 5096.357400000001
$ws.Range("M139").Value = 81.03999999999996
$ws.Range("N139").Value = -15376.3574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5236.8184
$ws.Range("I80").Value = 4650.625
$ws.Range("J80").Value = 6800
$ws.Range("K80").Value = 4650.625
$ws.Range("L80").Value = 6800
$ws.Range("M80").Value = -3652.625
$ws.Range("N80").Value = -8796

$ws.Range("H83").Value = 5236.8184
$ws.Range("I83").Value = 4650.625
$ws.Range("J83").Value = 6800
$ws.Range("K83").Value = 23253.125
$ws.Range("L83").Value = 34000
$ws.Range("M83").Value = -18261.125
$ws.Range("N83").Value = -43984

$ws.Range("H92").Value = 24909
$ws.Range("J92").Value = 24909
$ws.Range("L92").Value = 24909
$ws.Range("N92").Value = -28653

$ws.Range("H132").Value = 3363.85
$ws.Range("I132").Value = 3328.077
$ws.Range("J132").Value = 3430.2856
$ws.Range("K132").Value = 9984.231
$ws.Range("L132").Value = 10290.8568
$ws.Range("M132").Value = -7454.231
$ws.Range("N132").Value = -15350.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1871.7273
$ws.Range("I7").Value = 1799.1428
$ws.Range("J7").Value = 1998.75
$ws.Range("K7").Value = 1799.1428
$ws.Range("L7").Value = 1998.75
$ws.Range("M7").Value = -1687.1428
$ws.Range("N7").Value = -2222.75

$ws.Range("H22").Value = 989
$ws.Range("I22").Value = 916.3333
$ws.Range("J22").Value = 1025.3334
$ws.Range("K22").Value = 916.3333
$ws.Range("L22").Value = 1025.3334
$ws.Range("M22").Value = -621.3333
$ws.Range("N22").Value = -1615.3334

$ws.Range("H27").Value = 989
$ws.Range("I27").Value = 916.3333
$ws.Range("J27").Value = 1025.3334
$ws.Range("K27").Value = 916.3333
$ws.Range("L27").Value = 1025.3334
$ws.Range("M27").Value = -809.3333
$ws.Range("N27").Value = -1239.3334

$ws.Range("H126").Value = 1871.7273
$ws.Range("I126").Value = 1799.1428
$ws.Range("J126").Value = 1998.75
$ws.Range("K126").Value = 5397.428400000001
$ws.Range("L126").Value = 5996.25
$ws.Range("M126").Value = -2927.428400000001
$ws.Range("N126").Value = -10936.25

$ws.Range("H132").Value = 2524.3103
$ws.Range("I132").Value = 2138.7222
$ws.Range("J132").Value = 3155.2727
$ws.Range("K132").Value = 6416.1666
$ws.Range("L132").Value = 9465.8181
$ws.Range("M132").Value = -3886.1666
$ws.Range("N132").Value = -14525.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 26400
$ws.Range("J98").Value = 26400
$ws.Range("L98").Value = 26400
$ws.Range("N98").Value = -32390

$ws.Range("H126").Value = 47619972
$ws.Range("I126").Value = 58824324
$ws.Range("K126").Value = 176472972
$ws.Range("M126").Value = -176470502

$ws.Range("H132").Value = 1775.4584
$ws.Range("I132").Value = 1192.5
$ws.Range("J132").Value = 3524.3333
$ws.Range("K132").Value = 3577.5
$ws.Range("L132").Value = 10572.9999
$ws.Range("M132").Value = -1047.5
$ws.Range("N132").Value = -15632.9999
